# Apply updated cryptocurrency price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain numeric-looking strings (e.g. "396.00", "0.0950").
# Force those cells to Text format first so Excel does not coerce them into numbers
# and strip significant trailing/leading zeros, matching the source text values.
$textPriceRows = @(5,6,9,10,11,14,15,18,20,22,23,24,26,27,29,30,32,33,34,35,37,40,41,43,44,46,48,51)
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "56.943.59"
$ws.Range("E2").Value = "  +0.82%  "

# Row 3
$ws.Range("D3").Value = "3.244.60"
$ws.Range("E3").Value = "  +0.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "396.00"
$ws.Range("E5").Value = "  -1.25%  "

# Row 6
$ws.Range("D6").Value = "107.75"
$ws.Range("E6").Value = "  -3.35%  "

# Row 7
$ws.Range("E7").Value = "  +4.13%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  -1.64%  "

# Row 10
$ws.Range("D10").Value = "39.15"
$ws.Range("E10").Value = "  -1.42%  "

# Row 11
$ws.Range("D11").Value = "0.0950"
$ws.Range("E11").Value = "  +5.15%  "

# Row 12
$ws.Range("E12").Value = "  +1.90%  "

# Row 13
$ws.Range("D13").Value = "3.751.71"
$ws.Range("E13").Value = "  +0.13%  "

# Row 14
$ws.Range("D14").Value = "8.29"
$ws.Range("E14").Value = "  +2.11%  "

# Row 15
$ws.Range("D15").Value = "18.86"
$ws.Range("E15").Value = "  -2.02%  "

# Row 16
$ws.Range("D16").Value = "3.252.45"
$ws.Range("E16").Value = "  -0.05%  "

# Row 17
$ws.Range("E17").Value = "  -4.30%  "

# Row 18
$ws.Range("D18").Value = "10.92"
$ws.Range("E18").Value = "  +1.69%  "

# Row 19
$ws.Range("D19").Value = "56.736.52"
$ws.Range("E19").Value = "  +0.81%  "

# Row 20
$ws.Range("D20").Value = "3.34"
$ws.Range("E20").Value = "  -3.20%  "

# Row 21
$ws.Range("E21").Value = "  +5.99%  "

# Row 22
$ws.Range("D22").Value = "12.81"
$ws.Range("E22").Value = "  -3.04%  "

# Row 23
$ws.Range("D23").Value = "291.95"
$ws.Range("E23").Value = "  -1.12%  "

# Row 24
$ws.Range("D24").Value = "74.11"
$ws.Range("E24").Value = "  -0.86%  "

# Row 25
$ws.Range("E25").Value = "  -2.51%  "

# Row 26
$ws.Range("D26").Value = "7.98"
$ws.Range("E26").Value = "  -3.20%  "

# Row 27
$ws.Range("D27").Value = "28.02"
$ws.Range("E27").Value = "  -0.89%  "

# Row 28
$ws.Range("E28").Value = "  -0.31%  "

# Row 29
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "7.23"
$ws.Range("E29").Value = "  -3.84%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.169"
$ws.Range("E30").Value = "  -1.79%  "

# Row 31
$ws.Range("E31").Value = "  +0.07%  "

# Row 32
$ws.Range("D32").Value = "0.110"
$ws.Range("E32").Value = "  -2.65%  "

# Row 33
$ws.Range("D33").Value = "11.16"
$ws.Range("E33").Value = "  -2.04%  "

# Row 34
$ws.Range("D34").Value = "40.99"
$ws.Range("E34").Value = "  +11.23%  "

# Row 35
$ws.Range("D35").Value = "0.0483"
$ws.Range("E35").Value = "  -2.76%  "

# Row 36
$ws.Range("E36").Value = "  +0.85%  "

# Row 37
$ws.Range("D37").Value = "51.13"
$ws.Range("E37").Value = "  -0.73%  "

# Row 38
$ws.Range("E38").Value = "  +0.05%  "

# Row 40
$ws.Range("D40").Value = "2.97"
$ws.Range("E40").Value = "  -5.66%  "

# Row 41
$ws.Range("D41").Value = "137.31"
$ws.Range("E41").Value = "  +0.33%  "

# Row 42
$ws.Range("E42").Value = "  +1.59%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "1.87"
$ws.Range("E43").Value = "  -3.64%  "

# Row 44
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "3.91"
$ws.Range("E44").Value = "  -3.41%  "

# Row 45
$ws.Range("E45").Value = "  -1.27%  "

# Row 46
$ws.Range("D46").Value = "16.62"
$ws.Range("E46").Value = "  -3.68%  "

# Row 47
$ws.Range("E47").Value = "  +7.00%  "

# Row 48
$ws.Range("D48").Value = "22.27"
$ws.Range("E48").Value = "  -1.84%  "

# Row 49
$ws.Range("D49").Value = "2.154.10"
$ws.Range("E49").Value = "  -0.31%  "

# Row 50
$ws.Range("E50").Value = "  -4.17%  "

# Row 51
$ws.Range("D51").Value = "1.94"
$ws.Range("E51").Value = "  -8.74%  "
